$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.024.04"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.644.83"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.57"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "1.871.30"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "1.641.51"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.13"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "25.995.02"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.39"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +7.82%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.60"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0501"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.907"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "1.133.62"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.542"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.51"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.37"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "1.779.89"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.77"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.71"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  +0.29%  "
